$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the "primary goal" justification paragraph (5 local edits that
#    together turn the old sentence into the new one).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "use with stroke prediction (using dummy data), in a python-flask",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "use with a stroke prediction AI model trained on a provided dataset, in a python-flask",
    2) | Out-Null

$d.Content.Find.Execute(
    "achieve the desired outcomes, careful planning, and",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "achieve this desired outcome, careful planning and",
    2) | Out-Null

$d.Content.Find.Execute(
    "best practices and industry standards will need",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "best practices with industry standard techniques will need",
    2) | Out-Null

$apos = [char]8217
$d.Content.Find.Execute(
    "ethical considerations of data protection for patient${apos}s records",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ethical considerations for data protection of patient${apos}s records",
    2) | Out-Null

$d.Content.Find.Execute(
    "relevant regulatory standards.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "relevant regulatory requirements.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Capitalise a handful of tool / service names.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.MatchCase = $true
$rng.Find.Execute("bandit", $true, $false, $false, $false, $false, $true, 1, $false, "Bandit", 2) | Out-Null

$rng = $d.Content
$rng.Find.MatchCase = $true
$rng.Find.Execute("semgrep", $true, $false, $false, $false, $false, $true, 1, $false, "Semgrep", 2) | Out-Null

$rng = $d.Content
$rng.Find.MatchCase = $true
$rng.Find.Execute("github", $true, $false, $false, $false, $false, $true, 1, $false, "Github", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the "Design Methodology"
#    paragraph to the end of the SAST/testing-tools paragraph (right after
#    "...security mistakes.").
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a single-character placeholder right after the target sentence so we
# can anchor a bookmark to it, then shrink the bookmark to zero length by
# deleting the placeholder - this leaves a collapsed bookmark exactly at the
# desired position (a zero-length Range passed straight to Bookmarks.Add
# does not anchor reliably).
$anchorPhrase = "to check libraries for vulnerabilities or code for security mistakes."
$target = $d.Content
$target.Find.Execute($anchorPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.InsertAfter("X")

$markerRng = $d.Content
$markerRng.Find.Execute($anchorPhrase + "X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRng.MoveStart(1, $anchorPhrase.Length) | Out-Null

$d.Bookmarks.Add("_GoBack", $markerRng)
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$bmRange.Text = ""
